$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 2. Euclidean distance")
$ws.Columns.Item(1).Delete()
